$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Insert a fresh row at position 13. This shifts the old rows 13..23
#    down to 14..24, matching the new dimension A1:C24. Excel copies A12's
#    format down into the new A13, which we don't want (the new row has no
#    label in column A), so clear it right away.
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Clear()

# ---------------------------------------------------------------------------
# 2) Fill in the new row 13: the teacher's name, which used to (incorrectly)
#    sit next to "Objetivos:" on row 10, now gets its own row under
#    "Docentes responsáveis:" (row 12). Pull the column B/C formatting from
#    row 14 (the old row 13) so the new cells get the right wrap/color style.
# ---------------------------------------------------------------------------
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"

# ---------------------------------------------------------------------------
# 3) Correct / fill in the real text content for the various fields that
#    previously held placeholder / misplaced values.
# ---------------------------------------------------------------------------

# Objetivos: (row 10)
$objetivos = @'
Apresentar noções de mecânica dos fluidos e reologia, mediante estudo dos meios fluidos quando estáticos ou em movimento. Capacitar o aluno a modelar e resolver problemas de interesse em mecânica dos fluidos e reologia, com escolha adequada de hipóteses e aplicação de ferramentas correspondentes de solução.
'@
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Programa resumido: (row 14 after the insert)
$resumido = @'
Fundamentos de mecânica dos fluidos. Revisão de estática dos fluidos. Formulação integral e diferencial das equações de transporte de massa, energia e quantidade de movimento. Análise dimensional e semelhança. Escoamento incompressível de fluidos ideais e viscosos, regime laminar e turbulento. Equação de Navier-Stokes. Teoria da camada limite. Escoamento de fluidos não newtonianos. Formulação tensorial: tensão e deformação. Viscosidade e reometria. Viscoelasticidade. Aplicações.
'@
$ws.Range("B14").Value = $resumido
$ws.Range("C14").Value = $resumido

# Programa: (row 16 after the insert)
$programa = "Introdução: conceito de fluido; propriedades e conceito de contínuo; modelagem de processos de transferência; métodos de análise; dimensões e unidades.`nRevisão de estática de fluidos: equação básica da hidrostática, variação de pressão em um fluido estático; princípios de Stevin, de Pascal e de Arquimedes.`nFormulação integral das equações de transporte: teorema de transporte de Reynolds; aplicação para os princípios de conservação de massa, quantidade de movimento e energia; equação de Bernoulli.`nFormulação diferencial das equações de transporte: descrição do escoamento; forma diferencial: dos princípios de conservação de massa, quantidade de movimento e energia; formulação adimensional, análise dimensional e semelhança. Grupos adimensionais: número de Reynolds e número de Grashoff.`nEscoamento incompressível interno: equações de Euler; lei de Newton para a viscosidade, tensões de cisalhamento; equação de Navier-Stokes; regimes de escoamento: escoamento laminar e turbulento. Cálculo de perda de carga (distribuída e localizada), coeficiente de atrito. `nEscoamento incompressível externo: introdução à camada limite; escoamento ao redor de corpos, força da arraste.`nIntrodução a reologia. Definição e formulação tensorial de tensão e deformação. Tipos de deformação e escoamento de materiais. Equações fundamentais da reologia. Escoamento de fluidos newtonianos e não newtonianos. Viscosimetria e reometria. Reologia de sistemas dispersos. Colóides e emulsões. Soluções diluídas. Viscosimetria capilar. Aplicações."
$ws.Range("B16").Value = $programa
$ws.Range("C16").Value = $programa

# Método: (row 19 after the insert)
$metodo = @'
A avaliação será feita por meio de duas provas escritas P1 e P2 e por listas de exercícios e relatórios.
'@
$ws.Range("B19").Value = $metodo
$ws.Range("C19").Value = $metodo

# Critério: (row 20 after the insert)
$criterio = @'
A Nota final (NF) será calculada pela média ponderada das provas escritas e pela média dos trabalhos TR da seguinte maneira: NF = (P1 + 2*P2 + TR)/4
'@
$ws.Range("B20").Value = $criterio
$ws.Range("C20").Value = $criterio

# Norma de recuperação: (row 21 after the insert)
$recuperacao = @'
A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2
'@
$ws.Range("B21").Value = $recuperacao
$ws.Range("C21").Value = $recuperacao

# Bibliografia: (row 22 after the insert) - new multi-line content
$biblio = "BIRD,R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. LTC Editora, 2004.`nFOX, R. W., McDONALD, A. T. Introdução à Mecânica dos Fluidos. LTC Editora, 2001.`nSISSOM, L. E., PITTS, D. R. Fenômenos de Transporte. Ed. Guanabara, 1988.`nSCHRAMM, G. Reologia e Reometria. Editora Artliber, 2006.`nMANRICH, S.; PESSAN, L.A. Reologia: Conceitos Básicos, Editora UFSCar, 1987.`nMALKIN, A. Rheology Fundamentals. ChemTec Publishing, 1994."
$ws.Range("B22").Value = $biblio
$ws.Range("C22").Value = $biblio

# ---------------------------------------------------------------------------
# 4) Row height adjustments that accompany the content changes.
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).RowHeight = 60            # "Short syllabus:" row shrinks from 120 to 60
$ws.Rows.Item(21).RowHeight = 60            # "Norma de recuperação:" row shrinks from 120 to 60
$ws.Rows.Item(22).RowHeight = 120           # "Bibliografia:" row grows to 120 (has long text now)
